$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet1 ("MT") - sequence number (AK2) bumped from 1 to 3 (kept as text)
$ws1.Range("AK2").NumberFormat = "@"
$ws1.Range("AK2").Value = "3"

# Date (K2) updated on all sheets: 05-21-2024 -> 30-05-2024
$ws1.Range("K2").Value = "30-05-2024"
$ws2.Range("K2").Value = "30-05-2024"
$ws3.Range("K2").Value = "30-05-2024"
$ws4.Range("K2").Value = "30-05-2024"

# DateandTime (N2) + SequenceNo/ET id (AG2) on Sheet1
$ws1.Range("N2").Value = "30-05-2024 01:45:51 PM"
$ws1.Range("AG2").Value = "ET461"

# DateandTime (N2) on Sheet2, Sheet3, Sheet4
$ws2.Range("N2").Value = "30-05-2024 01:51:52 PM"
$ws3.Range("N2").Value = "30-05-2024 01:51:52 PM"
$ws4.Range("N2").Value = "30-05-2024 01:51:52 PM"

# SequenceNo / ET id (AG2) on Sheet2
$ws2.Range("AG2").Value = "ET462"

# SequenceNo / ET id (AG2) on Sheet3 and Sheet4
$ws3.Range("AG2").Value = "ET463"
$ws4.Range("AG2").Value = "ET463"
